$d = $word.ActiveDocument

# Build the OOXML fragment for the two new list-paragraphs (same list /
# style as the preceding "person" bullet, numId 7) plus a trailing blank
# paragraph, matching the target diff exactly.
$frag  = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr>'
$frag += '<w:r><w:t>SELECT * FROM `animal` WHERE type LIKE "%a%";</w:t></w:r></w:p>'
$frag += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr>'
$frag += '<w:r><w:t>SELECT * FROM `animal` WHERE type LIKE "%</w:t></w:r>'
$frag += '<w:proofErr w:type="spellStart"/><w:r><w:t>e</w:t></w:r><w:r><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$frag += '<w:r><w:t>%";</w:t></w:r></w:p>'
$frag += '<w:p/>'

$xml  = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>'
$xml += '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
$xml += '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">'
$xml += '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xml += $frag
$xml += '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Insert at the very end of the document's story so the existing last
# paragraph ("...person... LIKE "B%";") is left untouched and the new
# paragraphs are appended after it, before the final section break.
$r = $d.Content
$r.Collapse(0)
$r.InsertXML($xml) | Out-Null
